$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.073.58"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.012.37"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'226.17"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'0.598"
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'54.92"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("D12").Value = "2.309.87"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "'14.04"
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").Value = "'19.81"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "'0.737"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").Value = "2.016.22"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "36.985.98"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'6.26"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").Value = "'68.24"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "0.0₃0814"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "'222.74"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("D26").Value = "'164.54"
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  -6.23%  "
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").Value = "'18.56"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("E30").Value = "  -7.90%  "
$ws.Range("D31").Value = "'0.116"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").Value = "'4.49"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("D39").Value = "'5.34"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "1.456.46"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "'94.88"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").Value = "'2.77"
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("D44").Value = "'0.0905"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("E45").Value = "  -4.14%  "
$ws.Range("D46").Value = "'15.92"
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").Value = "'7.12"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "2.197.93"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "'3.58"
$ws.Range("E51").Value = "  -4.18%  "
